$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("Data")
$wsData.Rows("2:3").Insert()
$wsData.Rows("2:3").Delete()
Write-Output "ok"
